# Auto-generated Excel COM-interop script to apply market-price/profit
# corrections across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3757.25
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3757.25
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 11271.75
$ws.Range("N69").Value = -13019.75
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 3757.25
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3757.25
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 33815.25
$ws.Range("N72").Value = -42551.25
$ws.Range("M72").ClearContents()

$ws.Range("H135").Value = 100001570
$ws.Range("I135").Value = 35715960
$ws.Range("J135").Value = 1000000000
$ws.Range("K135").Value = 321443640
$ws.Range("L135").Value = 9000000000
$ws.Range("M135").Value = -321441105
$ws.Range("N135").Value = -9000005070

$ws.Range("H137").Value = 3427.2273
$ws.Range("I137").Value = 2524.9375
$ws.Range("J137").Value = 5833.3335
$ws.Range("K137").Value = 7574.8125
$ws.Range("L137").Value = 17500.0005
$ws.Range("M137").Value = -5024.8125
$ws.Range("N137").Value = -22600.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9211.948
$ws.Range("I61").Value = 5449.963
$ws.Range("J61").Value = 17676.416
$ws.Range("K61").Value = 5449.963
$ws.Range("L61").Value = 17676.416
$ws.Range("M61").Value = -5237.963
$ws.Range("N61").Value = -18100.416

$ws.Range("H74").Value = 3881.7874
$ws.Range("I74").Value = 1823.3684
$ws.Range("J74").Value = 12572.889
$ws.Range("K74").Value = 1823.3684
$ws.Range("L74").Value = 12572.889
$ws.Range("M74").Value = -949.3684000000001
$ws.Range("N74").Value = -14320.889

$ws.Range("H77").Value = 3881.7874
$ws.Range("I77").Value = 1823.3684
$ws.Range("J77").Value = 12572.889
$ws.Range("K77").Value = 9116.842000000001
$ws.Range("L77").Value = 62864.44499999999
$ws.Range("M77").Value = -4748.842000000001
$ws.Range("N77").Value = -71600.44499999999

$ws.Range("H132").Value = 2506.7742
$ws.Range("I132").Value = 1938
$ws.Range("J132").Value = 3701.2
$ws.Range("K132").Value = 5814
$ws.Range("L132").Value = 11103.6
$ws.Range("M132").Value = -3284
$ws.Range("N132").Value = -16163.6

$ws.Range("H136").Value = 9211.948
$ws.Range("I136").Value = 5449.963
$ws.Range("J136").Value = 17676.416
$ws.Range("K136").Value = 16349.889
$ws.Range("L136").Value = 53029.24800000001
$ws.Range("M136").Value = -13799.889
$ws.Range("N136").Value = -58129.24800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 28193.334
$ws.Range("J74").Value = 28193.334
$ws.Range("L74").Value = 28193.334
$ws.Range("N74").Value = -30065.334

$ws.Range("H77").Value = 28193.334
$ws.Range("J77").Value = 28193.334
$ws.Range("L77").Value = 84580.00199999999
$ws.Range("N77").Value = -93940.00199999999

$ws.Range("H94").Value = 1364.6
$ws.Range("I94").Value = 1364.6
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1364.6
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -913.5999999999999
$ws.Range("N94").ClearContents()

$ws.Range("H134").Value = 34617.453
$ws.Range("I134").Value = 2460.65
$ws.Range("J134").Value = 93084.37
$ws.Range("K134").Value = 7381.950000000001
$ws.Range("L134").Value = 279253.11
$ws.Range("M134").Value = -4846.950000000001
$ws.Range("N134").Value = -284323.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5576.636
$ws.Range("I31").Value = 9164.23
$ws.Range("J31").Value = 3244.7
$ws.Range("K31").Value = 9164.23
$ws.Range("L31").Value = 3244.7
$ws.Range("M31").Value = -8869.23
$ws.Range("N31").Value = -3834.7

$ws.Range("H34").Value = 5576.636
$ws.Range("I34").Value = 9164.23
$ws.Range("J34").Value = 3244.7
$ws.Range("K34").Value = 9164.23
$ws.Range("L34").Value = 3244.7
$ws.Range("M34").Value = -8962.23
$ws.Range("N34").Value = -3648.7

$ws.Range("H58").Value = 3032989.5
$ws.Range("I58").Value = 6061998.5
$ws.Range("J58").Value = 3980.0667
$ws.Range("K58").Value = 6061998.5
$ws.Range("L58").Value = 3980.0667
$ws.Range("M58").Value = -6061795.5
$ws.Range("N58").Value = -4386.066699999999

$ws.Range("H99").Value = 1355.7693
$ws.Range("I99").Value = 1271.1
$ws.Range("J99").Value = 1638
$ws.Range("K99").Value = 1271.1
$ws.Range("L99").Value = 1638
$ws.Range("M99").Value = 226.9000000000001
$ws.Range("N99").Value = -4634

$ws.Range("H126").Value = 1355.7693
$ws.Range("I126").Value = 1271.1
$ws.Range("J126").Value = 1638
$ws.Range("K126").Value = 3813.3
$ws.Range("L126").Value = 4914
$ws.Range("M126").Value = -1343.3
$ws.Range("N126").Value = -9854

$ws.Range("H132").Value = 2357.5095
$ws.Range("I132").Value = 2051.9412
$ws.Range("J132").Value = 2904.3157
$ws.Range("K132").Value = 6155.823600000001
$ws.Range("L132").Value = 8712.947100000001
$ws.Range("M132").Value = -3625.823600000001
$ws.Range("N132").Value = -13772.9471

$ws.Range("H134").Value = 2408.543
$ws.Range("I134").Value = 2151.2593
$ws.Range("J134").Value = 3276.875
$ws.Range("K134").Value = 6453.777900000001
$ws.Range("L134").Value = 9830.625
$ws.Range("M134").Value = -3918.777900000001
$ws.Range("N134").Value = -14900.625

$ws.Range("H136").Value = 3032989.5
$ws.Range("I136").Value = 6061998.5
$ws.Range("J136").Value = 3980.0667
$ws.Range("K136").Value = 18185995.5
$ws.Range("L136").Value = 11940.2001
$ws.Range("M136").Value = -18183445.5
$ws.Range("N136").Value = -17040.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1329.5428
$ws.Range("J131").Value = 1003.4667
$ws.Range("L131").Value = 3010.4001
$ws.Range("N131").Value = -13090.4001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4389.0625
$ws.Range("I80").Value = 2601
$ws.Range("J80").Value = 5201.8184
$ws.Range("K80").Value = 2601
$ws.Range("L80").Value = 5201.8184
$ws.Range("M80").Value = -1603
$ws.Range("N80").Value = -7197.8184

$ws.Range("H83").Value = 4389.0625
$ws.Range("I83").Value = 2601
$ws.Range("J83").Value = 5201.8184
$ws.Range("K83").Value = 13005
$ws.Range("L83").Value = 26009.092
$ws.Range("M83").Value = -8013
$ws.Range("N83").Value = -35993.092

$ws.Range("H132").Value = 6192.5864
$ws.Range("I132").Value = 2363
$ws.Range("J132").Value = 16245.25
$ws.Range("K132").Value = 7089
$ws.Range("L132").Value = 48735.75
$ws.Range("M132").Value = -4559
$ws.Range("N132").Value = -53795.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4055.4644
$ws.Range("I132").Value = 4097.1875
$ws.Range("J132").Value = 3999.8333
$ws.Range("K132").Value = 12291.5625
$ws.Range("L132").Value = 11999.4999
$ws.Range("M132").Value = -9761.5625
$ws.Range("N132").Value = -17059.4999

$ws.Range("H136").Value = 5872.4194
$ws.Range("I136").Value = 3130.7144
$ws.Range("J136").Value = 8130.294
$ws.Range("K136").Value = 9392.143199999999
$ws.Range("L136").Value = 24390.882
$ws.Range("M136").Value = -6842.143199999999
$ws.Range("N136").Value = -29490.882

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 36635.5
$ws.Range("J68").Value = 36635.5
$ws.Range("L68").Value = 36635.5
$ws.Range("N68").Value = -38257.5

$ws.Range("H71").Value = 36635.5
$ws.Range("J71").Value = 36635.5
$ws.Range("L71").Value = 109906.5
$ws.Range("N71").Value = -118018.5

$ws.Range("H132").Value = 2010.0541
$ws.Range("I132").Value = 1227.4584
$ws.Range("J132").Value = 3454.8462
$ws.Range("K132").Value = 3682.3752
$ws.Range("L132").Value = 10364.5386
$ws.Range("M132").Value = -1152.3752
$ws.Range("N132").Value = -15424.5386

$ws.Range("H136").Value = 6766.853
$ws.Range("I136").Value = 3018.625
$ws.Range("J136").Value = 10098.611
$ws.Range("K136").Value = 9055.875
$ws.Range("L136").Value = 30295.833
$ws.Range("M136").Value = -6505.875
$ws.Range("N136").Value = -35395.833
